# Refactor the format of the worksheet:
#  - Swap columns B and C (their content + column widths)
#  - Swap columns D and E (their content + column widths)
#  - Leave column D selected, matching the state Excel ends up
#    in after the user performs the column rearrangement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column C into column B's position (B moves right to C).
$ws.Columns("C").Cut()
$ws.Columns("B").Insert()

# Swap column E into column D's position (D moves right to E).
$ws.Columns("E").Cut()
$ws.Columns("D").Insert()

# Leave the selection on column D, as recorded in the saved workbook.
$null = $ws.Columns("D").Select()
